$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4: replace numeric values with strings "hello" / "hi"
$ws.Range("A4").Value = "hello"
$ws.Range("B4").Value = "hi"

# Add new row 5 with numeric values
$ws.Range("A5").Value = 23
$ws.Range("B5").Value = 32

# Update selection to A5, matching the diff
$ws.Range("A5").Select()
